$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  B=345; C=5.306968212127686},
    @{Row=3;  B=390; C=7.739184379577637},
    @{Row=4;  B=289; C=5.619730710983276},
    @{Row=5;  B=272; C=4.020173311233521},
    @{Row=6;  B=404; C=7.428812503814697},
    @{Row=7;  B=296; C=6.137551546096802},
    @{Row=8;  B=328; C=7.76811957359314},
    @{Row=9;  B=346; C=21.91056942939758},
    @{Row=10; B=379; C=28.65278673171997},
    @{Row=11; B=296; C=7.16151762008667},
    @{Row=12; B=480; C=14.11760377883911},
    @{Row=13; B=406; C=10.33298563957214},
    @{Row=14; B=685; C=12.52477717399597},
    @{Row=15; B=126; C=6.341109752655029},
    @{Row=16; B=113; C=4.343533277511597},
    @{Row=17; B=162; C=8.937321424484253},
    @{Row=18; B=169; C=10.12352561950684},
    @{Row=19; B=266; C=10.86786365509033},
    @{Row=20; B=67;  C=1.777767896652222},
    @{Row=21; B=127; C=6.24271035194397},
    @{Row=22; B=221; C=14.68635034561157},
    @{Row=23; B=311; C=25.20130062103271},
    @{Row=24; B=386; C=34.20357608795166}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
}
